# Apply crypto price/volume updates from the GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column (D) stored as text so values like "316.52" or
# "27.830.30" are not re-interpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.830.30'
$ws.Range("E2").Value = '  -2.82%  '
$ws.Range("D3").Value = '1.790.70'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '316.52'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.5344'
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("D8").Value = '0.3850'
$ws.Range("E8").Value = '  +2.12%  '
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("D10").Value = '41.30'
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("D11").Value = '1.083'
$ws.Range("E11").Value = '  -2.91%  '
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '6.170'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '7.417'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '20.27'
$ws.Range("E15").Value = '  -2.06%  '
$ws.Range("D16").Value = '1.790.65'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '88.08'
$ws.Range("E17").Value = '  -2.39%  '
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").Value = '0.06525'
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '17.18'
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").Value = '5.945'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").Value = '27.874.08'
$ws.Range("E23").Value = '  -2.70%  '
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").Value = '2.092'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").Value = '156.84'
$ws.Range("E26").Value = '  -1.58%  '
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("D28").Value = '1.997.68'
$ws.Range("D29").Value = '2.284'
$ws.Range("E29").Value = '  -2.79%  '
$ws.Range("D30").Value = '121.38'
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("D31").Value = '0.1094'
$ws.Range("E31").Value = '  +3.45%  '
$ws.Range("D32").Value = '1.094'
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").Value = '3.646'
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("D34").Value = '5.483'
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("D35").Value = '0.06954'
$ws.Range("E35").Value = '  +8.76%  '
$ws.Range("D36").Value = '0.2194'
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("D37").Value = '0.02259'
$ws.Range("E37").Value = '  -2.21%  '
$ws.Range("D38").Value = '5.050'
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").Value = '11.23'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '8.349'
$ws.Range("E40").Value = '  -4.69%  '
$ws.Range("D41").Value = '0.6081'
$ws.Range("E41").Value = '  -2.03%  '
$ws.Range("E42").Value = '  -3.19%  '
$ws.Range("D43").Value = '1.415'
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").Value = '13.24'
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("D45").Value = '3.679'
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("E46").Value = '  -2.84%  '
$ws.Range("D47").Value = '124.33'
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("E50").Value = '  -1.39%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '71.18'
$ws.Range("E51").Value = '  -1.34%  '
